$wb = $excel.ActiveWorkbook
$app = $wb.Application

$wsPending = $wb.Worksheets.Item("DPLKKPS131-002")
$wsVerif   = $wb.Worksheets.Item("DPLKKPS131-001")

# Update the rekening (account) number referenced in the "Lanjut Verifikasi" test
# case (sheet DPLKKPS131-002): both the free-text preparation notes (F2) and the
# dedicated NOMOR_REKENING column (P2) need the new account number.
$newPrepText = "Username : 30603;`n" + `
    "Password : bni1234;`n" + `
    "Role : Asisten Settlement;`n" + `
    "Keterangan Perubahan : KEP.TRX.436;`n" + `
    "Pilih Perusahaan : 000000029;`n" + `
    "No. Rekening : 805255179;`n" + `
    "Status Register : 1 : Lanjutkan ke Verifikasi;`n" + `
    "Keterangan Register : KEP.TRX.436 Lanjut Verifikasi"

$wsPending.Range("F2").Value = $newPrepText
$wsPending.Range("P2").Value = 805255179

# Restore view/selection: DPLKKPS131-002 is no longer the focused tab, and its
# selection resets back to the top-left corner of the data (A2).
$wsPending.Activate()
$win = $app.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$wsPending.Range("A2").Select() | Out-Null

# DPLKKPS131-001 becomes the active tab, with R2 selected.
$wsVerif.Activate()
$wsVerif.Range("R2").Select() | Out-Null
